$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B = New-Object 'object[,]' 24,3
$arr_B[0,0] = 0.7256381965642049
$arr_B[0,1] = 0.2458626454080388
$arr_B[0,2] = 0.1433199007234407
$arr_B[1,0] = 0.6848036183984902
$arr_B[1,1] = 0.2346601551900562
$arr_B[1,2] = 0.1379964880744922
$arr_B[2,0] = 0.6601359475813808
$arr_B[2,1] = 0.2279324937062484
$arr_B[2,2] = 0.1348052093490253
$arr_B[3,0] = 0.6501854823749511
$arr_B[3,1] = 0.2252286787129094
$arr_B[3,2] = 0.133524198356227
$arr_B[4,0] = 0.6485393642033443
$arr_B[4,1] = 0.2247819890743301
$arr_B[4,2] = 0.1333126634302459
$arr_B[5,0] = 0.6600013399113891
$arr_B[5,1] = 0.2278958763953653
$arr_B[5,2] = 0.1347878543640206
$arr_B[6,0] = 0.7114743691977594
$arr_B[6,1] = 0.2419686597786495
$arr_B[6,2] = 0.1414683496229685
$arr_B[7,0] = 0.8156326048110429
$arr_B[7,1] = 0.2707691083324164
$arr_B[7,2] = 0.1551823385378839
$arr_B[8,0] = 0.8941392799015944
$arr_B[8,1] = 0.2926761806544675
$arr_B[8,2] = 0.1656334743577474
$arr_B[9,0] = 0.930288971368725
$arr_B[9,1] = 0.3028076644183102
$arr_B[9,2] = 0.1704698985865463
$arr_B[10,0] = 0.944040894748241
$arr_B[10,1] = 0.3066682431483514
$arr_B[10,2] = 0.172313146116025
$arr_B[11,0] = 0.9410763774382076
$arr_B[11,1] = 0.3058357282554027
$arr_B[11,2] = 0.1719156452529376
$arr_B[12,0] = 0.9314190911554761
$arr_B[12,1] = 0.3031247945348525
$arr_B[12,2] = 0.1706213072869787
$arr_B[13,0] = 0.9255119066031625
$arr_B[13,1] = 0.3014674003139248
$arr_B[13,2] = 0.1698300240274051
$arr_B[14,0] = 0.8917856011935044
$arr_B[14,1] = 0.29201741500259
$arr_B[14,2] = 0.1653190544740539
$arr_B[15,0] = 0.8712074553343143
$arr_B[15,1] = 0.2862627434625438
$arr_B[15,2] = 0.1625727508076551
$arr_B[16,0] = 0.8594125435280944
$arr_B[16,1] = 0.2829684093292428
$arr_B[16,2] = 0.1610008827427549
$arr_B[17,0] = 0.855426045901055
$arr_B[17,1] = 0.2818556790884372
$arr_B[17,2] = 0.1604700038589613
$arr_B[18,0] = 0.873393780836409
$arr_B[18,1] = 0.2868737223250264
$arr_B[18,2] = 0.1628642992756966
$arr_B[19,0] = 0.9342539650092476
$arr_B[19,1] = 0.3039204087095015
$arr_B[19,2] = 0.1710011656219166
$arr_B[20,0] = 0.9743957959803424
$arr_B[20,1] = 0.3152014020503202
$arr_B[20,2] = 0.1763878436119768
$arr_B[21,0] = 0.9529378263209765
$arr_B[21,1] = 0.3091676611651053
$arr_B[21,2] = 0.1735065846770141
$arr_B[22,0] = 0.8724052317573694
$arr_B[22,1] = 0.2865974548637951
$arr_B[22,2] = 0.162732468351777
$arr_B[23,0] = 0.7871084251105742
$arr_B[23,1] = 0.262847562194338
$arr_B[23,2] = 0.1514065113894958
$ws.Range("B2:D25").Value = $arr_B

$arr_F = New-Object 'object[,]' 24,2
$arr_F[0,0] = 2.685327605700309
$arr_F[0,1] = 0.002537995726890081
$arr_F[1,0] = 2.670932625896057
$arr_F[1,1] = 0.002541995567351906
$arr_F[2,0] = 2.663430844192817
$arr_F[2,1] = 0.002544581224212066
$arr_F[3,0] = 2.660709535674869
$arr_F[3,1] = 0.002545667629590738
$arr_F[4,0] = 2.660277929566647
$arr_F[4,1] = 0.002545850006123088
$arr_F[5,0] = 2.663392784887861
$arr_F[5,1] = 0.002544595743016366
$arr_F[6,0] = 2.680086499664426
$arr_F[6,1] = 0.00253934800845867
$arr_F[7,0] = 2.723453815390769
$arr_F[7,1] = 0.002530081804465107
$arr_F[8,0] = 2.761837670276606
$arr_F[8,1] = 0.002523891772446916
$arr_F[9,0] = 2.780725087203351
$arr_F[9,1] = 0.002521208484050021
$arr_F[10,0] = 2.788083016373463
$arr_F[10,1] = 0.002520211348972834
$arr_F[11,0] = 2.786489198257811
$arr_F[11,1] = 0.002520425257930232
$arr_F[12,0] = 2.781326303213078
$arr_F[12,1] = 0.002521126069647035
$arr_F[13,0] = 2.778190682522805
$arr_F[13,1] = 0.002521557804183731
$arr_F[14,0] = 2.760632083670956
$arr_F[14,1] = 0.002524069791367323
$arr_F[15,0] = 2.750226183854721
$arr_F[15,1] = 0.002525644703059904
$arr_F[16,0] = 2.744375211828569
$arr_F[16,1] = 0.002526563036086876
$arr_F[17,0] = 2.742417209892537
$arr_F[17,1] = 0.002526876115394
$arr_F[18,0] = 2.75132001370541
$arr_F[18,1] = 0.002525475759425825
$arr_F[19,0] = 2.782837184197803
$arr_F[19,1] = 0.002520919710285674
$arr_F[20,0] = 2.804634570297594
$arr_F[20,1] = 0.002518052583456733
$arr_F[21,0] = 2.792890990988653
$arr_F[21,1] = 0.002519572743089306
$arr_F[22,0] = 2.75082508357761
$arr_F[22,1] = 0.002525552098537432
$arr_F[23,0] = 2.710579488138336
$arr_F[23,1] = 0.00253247957072507
$ws.Range("F2:G25").Value = $arr_F

$arr_I = New-Object 'object[,]' 24,3
$arr_I[0,0] = 1.298515895208098
$arr_I[0,1] = 0.3547417635487733
$arr_I[0,2] = 0.9630317485178921
$arr_I[1,0] = 1.298210087243255
$arr_I[1,1] = 0.3439979413373067
$arr_I[1,2] = 0.9118436157986594
$arr_I[2,0] = 1.29860705604721
$arr_I[2,1] = 0.3376178320811789
$arr_I[2,2] = 0.8809674899084996
$arr_I[3,0] = 1.298915685460784
$arr_I[3,1] = 0.335072185135644
$arr_I[3,2] = 0.8685242568205922
$arr_I[4,0] = 1.298975798146003
$arr_I[4,1] = 0.3346527568317867
$arr_I[4,2] = 0.8664664604011705
$arr_I[5,0] = 1.298610623933925
$arr_I[5,1] = 0.3375832809704207
$arr_I[5,2] = 0.8807991133694486
$arr_I[6,0] = 1.298288981386669
$arr_I[6,1] = 0.3509922187808741
$arr_I[6,2] = 0.94526706110085
$arr_I[7,0] = 1.302307281191361
$arr_I[7,1] = 0.3790158735812668
$arr_I[7,2] = 1.076096596410309
$arr_I[8,0] = 1.30810894253905
$arr_I[8,1] = 0.4006750154751018
$arr_I[8,2] = 1.174937477243617
$arr_I[9,0] = 1.311370474564427
$arr_I[9,1] = 0.4107643900057525
$arr_I[9,2] = 1.220501416417363
$arr_I[10,0] = 1.312695273350386
$arr_I[10,1] = 0.4146192324574827
$arr_I[10,2] = 1.237842088002481
$arr_I[11,0] = 1.312405959754486
$arr_I[11,1] = 0.4137874981607439
$arr_I[11,2] = 1.234103609224434
$arr_I[12,0] = 1.311477666914136
$arr_I[12,1] = 0.4110808433132576
$arr_I[12,2] = 1.22192630713036
$arr_I[13,0] = 1.310920753011501
$arr_I[13,1] = 0.409427400994204
$arr_I[13,2] = 1.21447864188687
$arr_I[14,0] = 1.307908336782191
$arr_I[14,1] = 0.400020427908629
$arr_I[14,2] = 1.171971872837901
$arr_I[15,0] = 1.306219877752582
$arr_I[15,1] = 0.3943102777772509
$arr_I[15,2] = 1.146049356314535
$arr_I[16,0] = 1.305307282700205
$arr_I[16,1] = 0.3910482007285765
$arr_I[16,2] = 1.131195940452784
$arr_I[17,0] = 1.305008344596352
$arr_I[17,1] = 0.389947532828046
$arr_I[17,2] = 1.126176531788275
$arr_I[18,0] = 1.306393554442593
$arr_I[18,1] = 0.3949158287219063
$arr_I[18,2] = 1.148803000290485
$arr_I[19,0] = 1.311747891936115
$arr_I[19,1] = 0.4118749232242038
$arr_I[19,2] = 1.225500723991615
$arr_I[20,0] = 1.315770373801129
$arr_I[20,1] = 0.4231582010143455
$arr_I[20,2] = 1.276131976516439
$arr_I[21,0] = 1.313575555243204
$arr_I[21,1] = 0.4171177784265865
$arr_I[21,2] = 1.249062857322343
$arr_I[22,0] = 1.306314854166658
$arr_I[22,1] = 0.3946419944827682
$arr_I[22,2] = 1.147557922793084
$arr_I[23,0] = 1.300720951290337
$arr_I[23,1] = 0.3712479836629825
$arr_I[23,2] = 1.040228428661578
$ws.Range("I2:K25").Value = $arr_I

$arr_N = New-Object 'object[,]' 24,1
$arr_N[0,0] = 2.444398860735767
$arr_N[1,0] = 2.461147951525128
$arr_N[2,0] = 2.472144976378189
$arr_N[3,0] = 2.476805450841823
$arr_N[4,0] = 2.477590130454843
$arr_N[5,0] = 2.472207104280898
$arr_N[6,0] = 2.450025843279249
$arr_N[7,0] = 2.412194509262704
$arr_N[8,0] = 2.387865558528446
$arr_N[9,0] = 2.37755280986579
$arr_N[10,0] = 2.373756395695736
$arr_N[11,0] = 2.374569177870598
$arr_N[12,0] = 2.377238294410219
$arr_N[13,0] = 2.378887382941059
$arr_N[14,0] = 2.388554730107501
$arr_N[15,0] = 2.394678811925488
$arr_N[16,0] = 2.398272231997254
$arr_N[17,0] = 2.399501089700792
$arr_N[18,0] = 2.394019541591348
$arr_N[19,0] = 2.376451355065441
$arr_N[20,0] = 2.365603835963171
$arr_N[21,0] = 2.371335225673221
$arr_N[22,0] = 2.39431737154537
$arr_N[23,0] = 2.42182112256183
$ws.Range("N2:N25").Value = $arr_N
